# Auto commit at 2025-11-09 10:04:41.97
# Append two new daily rows (2025-11-08 / serial 45969) for both stations
# to the bottom of the existing Sheet1 data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 - 四方坪站 (station "四方坪站" -> shared string index 4)
$ws.Range("A16").Value = 45969
$ws.Range("B16").Value = "四方坪站"
$ws.Range("C16").Value = 10999.84
$ws.Range("D16").Value = 9663.23
$ws.Range("E16").Value = 3708.77
$ws.Range("F16").Value = 462

# Row 17 - 高岭站 (station "高岭站" -> shared string index 5)
$ws.Range("A17").Value = 45969
$ws.Range("B17").Value = "高岭站"
$ws.Range("C17").Value = 4218.97
$ws.Range("D17").Value = 3625.19
$ws.Range("E17").Value = 1169.28
$ws.Range("F17").Value = 151

# Update the view state to match where the user ended up after typing the
# new rows in: scrolled down so row 7 is the first visible row, with the
# active selection sitting on M15.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("M15").Select() | Out-Null
